$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 173, pushing existing rows 173:234 down to 174:235
$ws.Rows.Item(173).Insert()

# Populate the newly inserted row 173 with the new weekly entry
$ws.Cells.Item(173, 1).Value = 10
$ws.Cells.Item(173, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(173, 3).Value = "La Araucanía"
$ws.Cells.Item(173, 4).Value = 44559
$ws.Cells.Item(173, 5).Value = 9
$ws.Cells.Item(173, 6).Value = 100112017
$ws.Cells.Item(173, 7).Value = "Apio"
$ws.Cells.Item(173, 8).Value = "Americana (o)"
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 155
$ws.Cells.Item(173, 11).Value = 8000
$ws.Cells.Item(173, 12).Value = 9000
$ws.Cells.Item(173, 13).Value = 8290
$ws.Cells.Item(173, 14).Value = "$/docena de matas"
$ws.Cells.Item(173, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(173, 16).Value = 1382
$ws.Cells.Item(173, 17).Value = 6
$ws.Cells.Item(173, 18).Value = "Hortaliza"

# Keep the date formatting consistent with the rest of column D
$ws.Cells.Item(173, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
